$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 237,2
$data[0,0] = 'param_TimeStep_starting_index'
$data[0,1] = 30.0
$data[1,0] = 'param_demand1_op_cost_starting_index'
$data[1,1] = 0
$data[2,0] = 'param_demand1_inv_cost_starting_index'
$data[2,1] = 0
$data[3,0] = 'param_demand2_inv_cost_starting_index'
$data[3,1] = 0
$data[4,0] = 'param_demand2_op_cost_starting_index'
$data[4,1] = 0
$data[5,0] = 'param_Q_net1_demand2_starting_index'
$data[5,1] = 0
$data[6,0] = 'param_net1_sell_thermal_starting_index'
$data[6,1] = 0
$data[7,0] = 'param_net1_buy_electric_starting_index'
$data[7,1] = 18.4
$data[8,0] = 'param_net1_sell_electric_starting_index'
$data[8,1] = 0
$data[9,0] = 'param_net1_emissions_starting_index'
$data[9,1] = 38.73842662151056
$data[10,0] = 'param_P_net1_bat2_starting_index'
$data[10,1] = 0
$data[11,0] = 'param_P_net1_heat_pump2_starting_index'
$data[11,1] = 20.0
$data[12,0] = 'param_P_net1_bat1_starting_index'
$data[12,1] = 0
$data[13,0] = 'param_P_net1_charging_station1_starting_index'
$data[13,1] = 0
$data[14,0] = 'param_P_to_net1_starting_index'
$data[14,1] = 0
$data[15,0] = 'param_P_net1_demand2_starting_index'
$data[15,1] = 0
$data[16,0] = 'param_P_from_net1_starting_index'
$data[16,1] = 40.0
$data[17,0] = 'param_P_net1_demand1_starting_index'
$data[17,1] = 0
$data[18,0] = 'param_Q_from_net1_starting_index'
$data[18,1] = 68.07677758962731
$data[19,0] = 'param_Q_to_net1_starting_index'
$data[19,1] = 0
$data[20,0] = 'param_net1_inv_cost_starting_index'
$data[20,1] = 0
$data[21,0] = 'param_P_net1_charging_station2_starting_index'
$data[21,1] = 0
$data[22,0] = 'param_P_net1_heat_pump1_starting_index'
$data[22,1] = 20.0
$data[23,0] = 'param_Q_net1_demand1_starting_index'
$data[23,1] = 68.07677758962731
$data[24,0] = 'param_net1_buy_thermal_starting_index'
$data[24,1] = 22.46533660457701
$data[25,0] = 'param_net2_buy_electric_starting_index'
$data[25,1] = 222.2871406905673
$data[26,0] = 'param_P_net2_bat1_starting_index'
$data[26,1] = 0
$data[27,0] = 'param_net2_inv_cost_starting_index'
$data[27,1] = 0
$data[28,0] = 'param_net2_buy_thermal_starting_index'
$data[28,1] = 167.448
$data[29,0] = 'param_net2_sell_thermal_starting_index'
$data[29,1] = 0
$data[30,0] = 'param_P_net2_bat2_starting_index'
$data[30,1] = 0
$data[31,0] = 'param_P_net2_charging_station1_starting_index'
$data[31,1] = 0
$data[32,0] = 'param_Q_from_net2_starting_index'
$data[32,1] = 837.24
$data[33,0] = 'param_P_net2_heat_pump2_starting_index'
$data[33,1] = 0
$data[34,0] = 'param_P_from_net2_starting_index'
$data[34,1] = 555.7178517264182
$data[35,0] = 'param_P_net2_demand2_starting_index'
$data[35,1] = 500.0
$data[36,0] = 'param_Q_net2_demand2_starting_index'
$data[36,1] = 837.24
$data[37,0] = 'param_net2_sell_electric_starting_index'
$data[37,1] = 0
$data[38,0] = 'param_P_to_net2_starting_index'
$data[38,1] = 0
$data[39,0] = 'param_P_net2_charging_station2_starting_index'
$data[39,1] = 0
$data[40,0] = 'param_Q_to_net2_starting_index'
$data[40,1] = 0
$data[41,0] = 'param_net2_emissions_starting_index'
$data[41,1] = 512.1395969667942
$data[42,0] = 'param_P_net2_heat_pump1_starting_index'
$data[42,1] = 0
$data[43,0] = 'param_P_net2_demand1_starting_index'
$data[43,1] = 55.71785172641819
$data[44,0] = 'param_Q_net2_demand1_starting_index'
$data[44,1] = 0
$data[45,0] = 'param_P_pv1_demand2_starting_index'
$data[45,1] = 0
$data[46,0] = 'param_pv1_op_cost_starting_index'
$data[46,1] = 1.0
$data[47,0] = 'param_pv1_emissions_starting_index'
$data[47,1] = 0
$data[48,0] = 'param_P_pv1_demand1_starting_index'
$data[48,1] = 0
$data[49,0] = 'param_P_pv1_bat1_starting_index'
$data[49,1] = 0
$data[50,0] = 'param_P_pv1_charging_station2_starting_index'
$data[50,1] = 0
$data[51,0] = 'param_P_pv1_bat2_starting_index'
$data[51,1] = 0
$data[52,0] = 'param_P_pv1_net2_starting_index'
$data[52,1] = 0
$data[53,0] = 'param_P_pv1_charging_station1_starting_index'
$data[53,1] = 0
$data[54,0] = 'param_P_from_pv1_starting_index'
$data[54,1] = 0
$data[55,0] = 'param_P_pv1_heat_pump2_starting_index'
$data[55,1] = 0
$data[56,0] = 'param_P_pv1_heat_pump1_starting_index'
$data[56,1] = 0
$data[57,0] = 'param_P_pv1_net1_starting_index'
$data[57,1] = 0
$data[58,0] = 'param_pv1_inv_cost_starting_index'
$data[58,1] = 0
$data[59,0] = 'param_P_pv2_charging_station1_starting_index'
$data[59,1] = 0
$data[60,0] = 'param_P_pv2_net1_starting_index'
$data[60,1] = 0
$data[61,0] = 'param_P_from_pv2_starting_index'
$data[61,1] = 1.8
$data[62,0] = 'param_P_pv2_bat1_starting_index'
$data[62,1] = 0
$data[63,0] = 'param_P_pv2_bat2_starting_index'
$data[63,1] = 0
$data[64,0] = 'param_P_pv2_demand2_starting_index'
$data[64,1] = 0
$data[65,0] = 'param_P_pv2_heat_pump1_starting_index'
$data[65,1] = 0
$data[66,0] = 'param_P_pv2_heat_pump2_starting_index'
$data[66,1] = 0
$data[67,0] = 'param_pv2_inv_cost_starting_index'
$data[67,1] = 0
$data[68,0] = 'param_P_pv2_charging_station2_starting_index'
$data[68,1] = 0
$data[69,0] = 'param_P_pv2_demand1_starting_index'
$data[69,1] = 1.8
$data[70,0] = 'param_P_pv2_net2_starting_index'
$data[70,1] = 0
$data[71,0] = 'param_pv2_op_cost_starting_index'
$data[71,1] = 1.0
$data[72,0] = 'param_pv2_emissions_starting_index'
$data[72,1] = 0.8999999999999999
$data[73,0] = 'param_bat1_K_ch_starting_index'
$data[73,1] = 0
$data[74,0] = 'param_P_bat1_net1_starting_index'
$data[74,1] = 0
$data[75,0] = 'param_P_bat1_heat_pump1_starting_index'
$data[75,1] = 0
$data[76,0] = 'param_P_bat1_demand2_starting_index'
$data[76,1] = 0
$data[77,0] = 'param_bat1_cumulated_aging_starting_index'
$data[77,1] = 0.000003166666666666667
$data[78,0] = 'param_bat1_op_cost_starting_index'
$data[78,1] = 1.0
$data[79,0] = 'param_bat1_emissions_starting_index'
$data[79,1] = 0
$data[80,0] = 'param_bat1_K_dis_starting_index'
$data[80,1] = 1.0
$data[81,0] = 'param_bat1_integer_starting_index'
$data[81,1] = 0
$data[82,0] = 'param_bat1_SOC_starting_index'
$data[82,1] = 0.3
$data[83,0] = 'param_P_bat1_charging_station1_starting_index'
$data[83,1] = 0
$data[84,0] = 'param_P_bat1_heat_pump2_starting_index'
$data[84,1] = 0
$data[85,0] = 'param_P_to_bat1_starting_index'
$data[85,1] = 0
$data[86,0] = 'param_bat1_inv_cost_starting_index'
$data[86,1] = 0
$data[87,0] = 'param_P_from_bat1_starting_index'
$data[87,1] = 0
$data[88,0] = 'param_P_bat1_demand1_starting_index'
$data[88,1] = 0
$data[89,0] = 'param_P_bat1_charging_station2_starting_index'
$data[89,1] = 0
$data[90,0] = 'param_P_bat1_net2_starting_index'
$data[90,1] = 0
$data[91,0] = 'param_bat1_SOC_max_starting_index'
$data[91,1] = 0.9999968333333333
$data[92,0] = 'param_P_bat2_net2_starting_index'
$data[92,1] = 0
$data[93,0] = 'param_P_bat2_demand1_starting_index'
$data[93,1] = 0
$data[94,0] = 'param_bat2_SOC_starting_index'
$data[94,1] = 0.5
$data[95,0] = 'param_P_bat2_charging_station1_starting_index'
$data[95,1] = 0
$data[96,0] = 'param_bat2_op_cost_starting_index'
$data[96,1] = 1.0
$data[97,0] = 'param_P_bat2_net1_starting_index'
$data[97,1] = 0
$data[98,0] = 'param_bat2_inv_cost_starting_index'
$data[98,1] = 0
$data[99,0] = 'param_bat2_emissions_starting_index'
$data[99,1] = 0
$data[100,0] = 'param_bat2_cumulated_aging_starting_index'
$data[100,1] = 0
$data[101,0] = 'param_P_bat2_demand2_starting_index'
$data[101,1] = 0
$data[102,0] = 'param_P_to_bat2_starting_index'
$data[102,1] = 0
$data[103,0] = 'param_P_bat2_charging_station2_starting_index'
$data[103,1] = 0
$data[104,0] = 'param_P_bat2_heat_pump1_starting_index'
$data[104,1] = 0
$data[105,0] = 'param_P_bat2_heat_pump2_starting_index'
$data[105,1] = 0
$data[106,0] = 'param_P_from_bat2_starting_index'
$data[106,1] = 0
$data[107,0] = 'param_bat2_K_dis_starting_index'
$data[107,1] = 1.0
$data[108,0] = 'param_bat2_K_ch_starting_index'
$data[108,1] = 0
$data[109,0] = 'param_bat2_SOC_max_starting_index'
$data[109,1] = 1.0
$data[110,0] = 'param_bat2_integer_starting_index'
$data[110,1] = 0
$data[111,0] = 'param_Q_CHP1_demand1_starting_index'
$data[111,1] = 40.0
$data[112,0] = 'param_P_CHP1_bat2_starting_index'
$data[112,1] = 0
$data[113,0] = 'param_P_CHP1_charging_station2_starting_index'
$data[113,1] = 0
$data[114,0] = 'param_P_CHP1_demand2_starting_index'
$data[114,1] = 0
$data[115,0] = 'param_CHP1_emissions_starting_index'
$data[115,1] = 4.83
$data[116,0] = 'param_P_CHP1_heat_pump1_starting_index'
$data[116,1] = 0
$data[117,0] = 'param_P_CHP1_bat1_starting_index'
$data[117,1] = 0
$data[118,0] = 'param_P_CHP1_net2_starting_index'
$data[118,1] = 0
$data[119,0] = 'param_CHP1_inv_cost_starting_index'
$data[119,1] = 0
$data[120,0] = 'param_P_CHP1_net1_starting_index'
$data[120,1] = 0
$data[121,0] = 'param_Q_CHP1_net2_starting_index'
$data[121,1] = 0
$data[122,0] = 'param_P_CHP1_charging_station1_starting_index'
$data[122,1] = 0
$data[123,0] = 'param_P_CHP1_heat_pump2_starting_index'
$data[123,1] = 0
$data[124,0] = 'param_P_from_CHP1_starting_index'
$data[124,1] = 20.0
$data[125,0] = 'param_Q_CHP1_net1_starting_index'
$data[125,1] = 0
$data[126,0] = 'param_P_CHP1_demand1_starting_index'
$data[126,1] = 20.0
$data[127,0] = 'param_CHP1_fuel_cons_starting_index'
$data[127,1] = 2.1
$data[128,0] = 'param_CHP1_op_cost_starting_index'
$data[128,1] = 10.5
$data[129,0] = 'param_Q_from_CHP1_starting_index'
$data[129,1] = 40.0
$data[130,0] = 'param_Q_CHP1_demand2_starting_index'
$data[130,1] = 0
$data[131,0] = 'param_P_CHP2_net2_starting_index'
$data[131,1] = 0
$data[132,0] = 'param_P_CHP2_bat2_starting_index'
$data[132,1] = 0
$data[133,0] = 'param_P_from_CHP2_starting_index'
$data[133,1] = 20.0
$data[134,0] = 'param_CHP2_op_cost_starting_index'
$data[134,1] = 10.5
$data[135,0] = 'param_Q_CHP2_net1_starting_index'
$data[135,1] = 0
$data[136,0] = 'param_P_CHP2_heat_pump2_starting_index'
$data[136,1] = 0
$data[137,0] = 'param_Q_CHP2_net2_starting_index'
$data[137,1] = 0
$data[138,0] = 'param_Q_CHP2_demand1_starting_index'
$data[138,1] = 40.0
$data[139,0] = 'param_CHP2_fuel_cons_starting_index'
$data[139,1] = 2.1
$data[140,0] = 'param_P_CHP2_charging_station1_starting_index'
$data[140,1] = 0
$data[141,0] = 'param_CHP2_inv_cost_starting_index'
$data[141,1] = 0
$data[142,0] = 'param_P_CHP2_bat1_starting_index'
$data[142,1] = 0
$data[143,0] = 'param_P_CHP2_demand1_starting_index'
$data[143,1] = 20.0
$data[144,0] = 'param_Q_from_CHP2_starting_index'
$data[144,1] = 40.0
$data[145,0] = 'param_Q_CHP2_demand2_starting_index'
$data[145,1] = 0
$data[146,0] = 'param_P_CHP2_heat_pump1_starting_index'
$data[146,1] = 0
$data[147,0] = 'param_P_CHP2_net1_starting_index'
$data[147,1] = 0
$data[148,0] = 'param_CHP2_emissions_starting_index'
$data[148,1] = 4.83
$data[149,0] = 'param_P_CHP2_demand2_starting_index'
$data[149,1] = 0
$data[150,0] = 'param_P_CHP2_charging_station2_starting_index'
$data[150,1] = 0
$data[151,0] = 'param_Q_solar_th1_net2_starting_index'
$data[151,1] = 0
$data[152,0] = 'param_Q_from_solar_th1_starting_index'
$data[152,1] = 0
$data[153,0] = 'param_Q_solar_th1_demand1_starting_index'
$data[153,1] = 0
$data[154,0] = 'param_Q_solar_th1_net1_starting_index'
$data[154,1] = 0
$data[155,0] = 'param_solar_th1_op_cost_starting_index'
$data[155,1] = 1.0
$data[156,0] = 'param_solar_th1_inv_cost_starting_index'
$data[156,1] = 0
$data[157,0] = 'param_solar_th1_emissions_starting_index'
$data[157,1] = 0
$data[158,0] = 'param_Q_solar_th1_demand2_starting_index'
$data[158,1] = 0
$data[159,0] = 'param_Q_solar_th2_net1_starting_index'
$data[159,1] = 0
$data[160,0] = 'param_Q_solar_th2_net2_starting_index'
$data[160,1] = 0
$data[161,0] = 'param_Q_solar_th2_demand2_starting_index'
$data[161,1] = 1.2
$data[162,0] = 'param_solar_th2_emissions_starting_index'
$data[162,1] = 0.6000000000000001
$data[163,0] = 'param_Q_from_solar_th2_starting_index'
$data[163,1] = 1.2
$data[164,0] = 'param_solar_th2_op_cost_starting_index'
$data[164,1] = 1.0
$data[165,0] = 'param_Q_solar_th2_demand1_starting_index'
$data[165,1] = 0
$data[166,0] = 'param_solar_th2_inv_cost_starting_index'
$data[166,1] = 0
$data[167,0] = 'param_P_pvt1_bat2_starting_index'
$data[167,1] = 0
$data[168,0] = 'param_P_from_pvt1_starting_index'
$data[168,1] = 0
$data[169,0] = 'param_pvt1_emissions_starting_index'
$data[169,1] = 0
$data[170,0] = 'param_pvt1_inv_cost_starting_index'
$data[170,1] = 0
$data[171,0] = 'param_P_pvt1_bat1_starting_index'
$data[171,1] = 0
$data[172,0] = 'param_P_pvt1_net2_starting_index'
$data[172,1] = 0
$data[173,0] = 'param_pvt1_op_cost_starting_index'
$data[173,1] = 1.0
$data[174,0] = 'param_P_pvt1_heat_pump2_starting_index'
$data[174,1] = 0
$data[175,0] = 'param_P_pvt1_charging_station2_starting_index'
$data[175,1] = 0
$data[176,0] = 'param_Q_pvt1_net1_starting_index'
$data[176,1] = 0
$data[177,0] = 'param_Q_from_pvt1_starting_index'
$data[177,1] = 0
$data[178,0] = 'param_P_pvt1_net1_starting_index'
$data[178,1] = 0
$data[179,0] = 'param_Q_pvt1_demand2_starting_index'
$data[179,1] = 0
$data[180,0] = 'param_P_pvt1_demand1_starting_index'
$data[180,1] = 0
$data[181,0] = 'param_P_pvt1_heat_pump1_starting_index'
$data[181,1] = 0
$data[182,0] = 'param_Q_pvt1_net2_starting_index'
$data[182,1] = 0
$data[183,0] = 'param_P_pvt1_charging_station1_starting_index'
$data[183,1] = 0
$data[184,0] = 'param_Q_pvt1_demand1_starting_index'
$data[184,1] = 0
$data[185,0] = 'param_P_pvt1_demand2_starting_index'
$data[185,1] = 0
$data[186,0] = 'param_P_pvt2_bat1_starting_index'
$data[186,1] = 0
$data[187,0] = 'param_P_pvt2_demand2_starting_index'
$data[187,1] = 0
$data[188,0] = 'param_P_pvt2_net1_starting_index'
$data[188,1] = 0
$data[189,0] = 'param_Q_from_pvt2_starting_index'
$data[189,1] = 1.56
$data[190,0] = 'param_P_from_pvt2_starting_index'
$data[190,1] = 1.2
$data[191,0] = 'param_P_pvt2_charging_station2_starting_index'
$data[191,1] = 0
$data[192,0] = 'param_P_pvt2_heat_pump2_starting_index'
$data[192,1] = 0
$data[193,0] = 'param_P_pvt2_bat2_starting_index'
$data[193,1] = 0
$data[194,0] = 'param_pvt2_emissions_starting_index'
$data[194,1] = 0.7800000000000001
$data[195,0] = 'param_Q_pvt2_demand1_starting_index'
$data[195,1] = 0
$data[196,0] = 'param_P_pvt2_charging_station1_starting_index'
$data[196,1] = 0
$data[197,0] = 'param_P_pvt2_demand1_starting_index'
$data[197,1] = 1.2
$data[198,0] = 'param_Q_pvt2_net2_starting_index'
$data[198,1] = 0
$data[199,0] = 'param_P_pvt2_net2_starting_index'
$data[199,1] = 0
$data[200,0] = 'param_P_pvt2_heat_pump1_starting_index'
$data[200,1] = 0
$data[201,0] = 'param_Q_pvt2_net1_starting_index'
$data[201,1] = 0
$data[202,0] = 'param_pvt2_inv_cost_starting_index'
$data[202,1] = 0
$data[203,0] = 'param_Q_pvt2_demand2_starting_index'
$data[203,1] = 1.56
$data[204,0] = 'param_pvt2_op_cost_starting_index'
$data[204,1] = 1.0
$data[205,0] = 'param_charging_station1_inv_cost_starting_index'
$data[205,1] = 0
$data[206,0] = 'param_charging_station1_op_cost_starting_index'
$data[206,1] = 0
$data[207,0] = 'param_charging_station1_emissions_starting_index'
$data[207,1] = 0
$data[208,0] = 'param_charging_station2_emissions_starting_index'
$data[208,1] = 0
$data[209,0] = 'param_charging_station2_inv_cost_starting_index'
$data[209,1] = 0
$data[210,0] = 'param_charging_station2_op_cost_starting_index'
$data[210,1] = 0
$data[211,0] = 'param_Q_heat_pump1_net1_starting_index'
$data[211,1] = 0
$data[212,0] = 'param_heat_pump1_op_cost_starting_index'
$data[212,1] = 8.561643835616438
$data[213,0] = 'param_heat_pump1_emissions_starting_index'
$data[213,1] = 2.76
$data[214,0] = 'param_Q_from_heat_pump1_starting_index'
$data[214,1] = 80.0
$data[215,0] = 'param_Q_to_heat_pump1_starting_index'
$data[215,1] = 0
$data[216,0] = 'param_P_from_heat_pump1_starting_index'
$data[216,1] = 0
$data[217,0] = 'param_heat_pump1_inv_cost_starting_index'
$data[217,1] = 0
$data[218,0] = 'param_Q_heat_pump1_net2_starting_index'
$data[218,1] = 0
$data[219,0] = 'param_Q_heat_pump1_demand1_starting_index'
$data[219,1] = 0
$data[220,0] = 'param_Q_heat_pump1_demand2_starting_index'
$data[220,1] = 80.0
$data[221,0] = 'param_P_to_heat_pump1_starting_index'
$data[221,1] = 20.0
$data[222,0] = 'param_P_to_heat_pump2_starting_index'
$data[222,1] = 20.0
$data[223,0] = 'param_Q_from_heat_pump2_starting_index'
$data[223,1] = 80.0
$data[224,0] = 'param_heat_pump2_emissions_starting_index'
$data[224,1] = 2.76
$data[225,0] = 'param_P_from_heat_pump2_starting_index'
$data[225,1] = 0
$data[226,0] = 'param_Q_heat_pump2_demand2_starting_index'
$data[226,1] = 80.0
$data[227,0] = 'param_Q_heat_pump2_net1_starting_index'
$data[227,1] = 0
$data[228,0] = 'param_Q_heat_pump2_net2_starting_index'
$data[228,1] = 0
$data[229,0] = 'param_Q_heat_pump2_demand1_starting_index'
$data[229,1] = 0
$data[230,0] = 'param_heat_pump2_op_cost_starting_index'
$data[230,1] = 8.561643835616438
$data[231,0] = 'param_heat_pump2_inv_cost_starting_index'
$data[231,1] = 0
$data[232,0] = 'param_Q_to_heat_pump2_starting_index'
$data[232,1] = 0
$data[233,0] = 'param_total_emissions_starting_index'
$data[233,1] = 566.9580235883047
$data[234,0] = 'param_total_sell_starting_index'
$data[234,1] = 0
$data[235,0] = 'param_total_buy_starting_index'
$data[235,1] = 430.6004772951443
$data[236,0] = 'param_total_operation_cost_starting_index'
$data[236,1] = 42.12328767123287
$ws.Range("A2:B238").Value = $data
